$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.40"
$ws.Range("E2").Value = "'6.27%"
$ws.Range("D3").Value = "'49.17"
$ws.Range("E3").Value = "'11.37%"
$ws.Range("D4").Value = "'5.339"
$ws.Range("E4").Value = "'4.94%"
$ws.Range("D5").Value = "'0.08067"
$ws.Range("E5").Value = "'4.77%"
$ws.Range("D6").Value = "'4.595"
$ws.Range("E6").Value = "'4.00%"
$ws.Range("D7").Value = "'1.351"
$ws.Range("E7").Value = "'29.05%"
$ws.Range("D8").Value = "'1.644"
$ws.Range("E8").Value = "'1.75%"
$ws.Range("D9").Value = "'0.1280"
$ws.Range("E9").Value = "'0.82%"
$ws.Range("D10").Value = "'0.1970"
$ws.Range("E10").Value = "'5.66%"
$ws.Range("D11").Value = "'0.09663"
$ws.Range("E11").Value = "'4.61%"
$ws.Range("D12").Value = "'0.04684"
$ws.Range("E12").Value = "'11.66%"
$ws.Range("D14").Value = "'0.001317"
$ws.Range("E14").Value = "'2.66%"
$ws.Range("D15").Value = "'0.04197"
$ws.Range("E15").Value = "'0.00%"
$ws.Range("D16").Value = "'0.005831"
$ws.Range("E16").Value = "'1.23%"
$ws.Range("D17").Value = "'3.342"
$ws.Range("E17").Value = "'-0.11%"
$ws.Range("D18").Value = "'2.446"
$ws.Range("E18").Value = "'4.97%"
$ws.Range("E19").Value = "'4.44%"
$ws.Range("D20").Value = "'8.034"
$ws.Range("E20").Value = "'-0.42%"
$ws.Range("D21").Value = "'0.1364"
$ws.Range("E21").Value = "'-2.48%"
$ws.Range("D22").Value = "'0.3092"
$ws.Range("E22").Value = "'-2.69%"
$ws.Range("D23").Value = "'0.001312"
$ws.Range("E23").Value = "'2.32%"
$ws.Range("D24").Value = "'0.004272"
$ws.Range("E24").Value = "'-3.15%"
$ws.Range("D25").Value = "'0.0001350"
$ws.Range("E25").Value = "'0.07%"
$ws.Range("D26").Value = "'0.0003540"
$ws.Range("D38").Value = "'0.02724"
$ws.Range("E38").Value = "'9.31%"
$ws.Range("D39").Value = "'0.06118"
$ws.Range("E39").Value = "'15.19%"
$ws.Range("D40").Value = "'0.01087"
$ws.Range("E40").Value = "'83.29%"
$ws.Range("D41").Value = "'0.008030"
$ws.Range("E41").Value = "'3.98%"
$ws.Range("D42").Value = "'0.1466"
$ws.Range("E42").Value = "'8.68%"
$ws.Range("D43").Value = "'0.007903"
$ws.Range("E43").Value = "'7.63%"
$ws.Range("D44").Value = "'0.008660"
$ws.Range("E44").Value = "'14.73%"
$ws.Range("D45").Value = "'0.3491"
$ws.Range("E45").Value = "'16.08%"
$ws.Range("D46").Value = "'0.00006856"
$ws.Range("E46").Value = "'3.09%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'0.05975"
$ws.Range("E48").Value = "'39.00%"
$ws.Range("D49").Value = "'0.004001"
$ws.Range("E49").Value = "'-4.77%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.06%"
